$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{D010F14D-2528-4D93-A074-09B1BEE4870D}")
